$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pin Map")
$ws.Columns("K").Insert()
